# Commit: "Use Then instead of Assert to match bdd syntax"
# The workbook uses BDD-style keywords in column A (Specification, Given a,
# When, ... Assert). Rename the "Assert" keyword to "Then" so the sheet
# matches standard Given/When/Then BDD terminology.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell that held the literal "Assert" keyword.
$ws.Range("A12").Value = "Then"

# The conditional-formatting rule that highlights the "Assert" keyword row
# needs to key off the new text too.
$rng = $ws.Range("A1:XFD1048576")
$fcs = $rng.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.Formula1 -eq '="Assert"') {
        $fc.Formula1 = '="Then"'
    }
}
